# Generate Report for Handback
#
# Renames the two handed-back file identifiers throughout the workbook:
#   c8766d38-3bf5-4885-b140-3c9fcaf444ec  ->  9993c348-c562-422b-8d38-0d8a9c505173
#   eca71562-2752-469a-8e40-58d52f15d7d0  ->  ffffa1923a7b-fa80-43ac-9549-8f327d1787e2
# and consolidates both locales onto a single new handoff/handback archive hash
# (9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5), with refreshed handoff/handback
# timestamps.

$GUID1 = "9993c348-c562-422b-8d38-0d8a9c505173"
$GUID2 = "ffffa1923a7b-fa80-43ac-9549-8f327d1787e2"
$HASH  = "9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5"

$zhcnXlf = "$GUID1.$HASH.zh-cn.xlf"
$dedeXlf = "$GUID1.$HASH.de-de.xlf"

$wb = $excel.ActiveWorkbook

# Note: a cell's displayed text and its hyperlink's displayed text are
# stored independently by this engine, so both the cell .Value and the
# matching Hyperlink.TextToDisplay need to be set explicitly below.

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$GUID1.md"
$wsOverview.Range("A3").Value = "$GUID2.md"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.TextToDisplay = "$GUID1.md"
    } elseif ($addr -eq "`$A`$3") {
        $hl.TextToDisplay = "$GUID2.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$GUID1.md"
$wsZh.Range("D2").Value = $zhcnXlf
$wsZh.Range("E2").Value = "2016-03-22 07:08:53"
$wsZh.Range("F2").Value = "$GUID1.md"
$wsZh.Range("G2").Value = $zhcnXlf
$wsZh.Range("H2").Value = "2016-03-22 07:09:15"

$wsZh.Range("A3").Value = "$GUID2.md"
$wsZh.Range("D3").Value = $zhcnXlf
$wsZh.Range("E3").Value = "2016-03-22 07:08:53"
$wsZh.Range("F3").Value = "$GUID2.md"
$wsZh.Range("G3").Value = $zhcnXlf
$wsZh.Range("H3").Value = "2016-03-22 07:09:15"

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.TextToDisplay = "$GUID1.md"
    } elseif ($addr -eq "`$D`$2") {
        $hl.TextToDisplay = $zhcnXlf
    } elseif ($addr -eq "`$F`$2") {
        $hl.TextToDisplay = "$GUID1.md"
    } elseif ($addr -eq "`$G`$2") {
        $hl.TextToDisplay = $zhcnXlf
    } elseif ($addr -eq "`$A`$3") {
        $hl.TextToDisplay = "$GUID2.md"
    } elseif ($addr -eq "`$D`$3") {
        $hl.TextToDisplay = $zhcnXlf
    } elseif ($addr -eq "`$F`$3") {
        $hl.TextToDisplay = "$GUID2.md"
    } elseif ($addr -eq "`$G`$3") {
        $hl.TextToDisplay = $zhcnXlf
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$GUID1.md"
$wsDe.Range("D2").Value = $dedeXlf
$wsDe.Range("E2").Value = "2016-03-22 07:08:57"
$wsDe.Range("F2").Value = "$GUID1.md"
$wsDe.Range("G2").Value = $dedeXlf
$wsDe.Range("H2").Value = "2016-03-22 07:09:21"

$wsDe.Range("A3").Value = "$GUID2.md"
$wsDe.Range("D3").Value = $dedeXlf
$wsDe.Range("E3").Value = "2016-03-22 07:08:57"
$wsDe.Range("F3").Value = "$GUID2.md"
$wsDe.Range("G3").Value = $dedeXlf
$wsDe.Range("H3").Value = "2016-03-22 07:09:21"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.TextToDisplay = "$GUID1.md"
    } elseif ($addr -eq "`$D`$2") {
        $hl.TextToDisplay = $dedeXlf
    } elseif ($addr -eq "`$F`$2") {
        $hl.TextToDisplay = "$GUID1.md"
    } elseif ($addr -eq "`$G`$2") {
        $hl.TextToDisplay = $dedeXlf
    } elseif ($addr -eq "`$A`$3") {
        $hl.TextToDisplay = "$GUID2.md"
    } elseif ($addr -eq "`$D`$3") {
        $hl.TextToDisplay = $dedeXlf
    } elseif ($addr -eq "`$F`$3") {
        $hl.TextToDisplay = "$GUID2.md"
    } elseif ($addr -eq "`$G`$3") {
        $hl.TextToDisplay = $dedeXlf
    }
}
